$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (closest achievable via ColumnWidth quantization)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666668
$ws.Columns.Item(2).ColumnWidth = 13.833333333333332

# Cell values A1:B32
$ws.Range("A1").Value = -0.38057104854618728
$ws.Range("B1").Value = 0.37927322379067618
$ws.Range("A2").Value = -0.26325215972463667
$ws.Range("B2").Value = 0.25962011509415106
$ws.Range("A3").Value = -0.09310364038294594
$ws.Range("B3").Value = 0.092647491785431413
$ws.Range("A4").Value = -0.15464092458532974
$ws.Range("B4").Value = 0.15382707976412568
$ws.Range("A5").Value = -0.14782708057912597
$ws.Range("B5").Value = 0.14619864114066949
$ws.Range("A6").Value = -0.053076732483280153
$ws.Range("B6").Value = 0.053035602033046025
$ws.Range("A7").Value = -0.033035603009425429
$ws.Range("B7").Value = 0.032978535090117944
$ws.Range("A8").Value = -0.012978536070807678
$ws.Range("B8").Value = 0.012967093156020404
$ws.Range("A9").Value = -0.0069670940036239415
$ws.Range("B9").Value = 0.0069640641355466215
$ws.Range("A10").Value = -0.00096406498484213898
$ws.Range("B10").Value = 0.00096617581801439201
$ws.Range("A11").Value = -0.051479374529051825
$ws.Range("B11").Value = 0.051401115218475013
$ws.Range("A12").Value = -0.045401116070890257
$ws.Range("B12").Value = 0.045152647270962909
$ws.Range("A13").Value = -0.039152648137294577
$ws.Range("B13").Value = 0.03908528541702605
$ws.Range("A14").Value = -0.027085286345517545
$ws.Range("B14").Value = 0.027053092543551038
$ws.Range("A15").Value = -0.021053093416464996
$ws.Range("B15").Value = 0.021027786094761858
$ws.Range("A16").Value = -0.015027786970258195
$ws.Range("B16").Value = 0.015004167906216281
$ws.Range("A17").Value = -0.0090041687852551178
$ws.Range("B17").Value = 0.0089999990911904248
$ws.Range("A18").Value = -0.10602106488071783
$ws.Range("B18").Value = 0.10586997962749223
$ws.Range("A19").Value = -0.027096764591751477
$ws.Range("B19").Value = 0.027013291199960321
$ws.Range("A20").Value = -0.01801329201488322
$ws.Range("B20").Value = 0.01800426176590264
$ws.Range("A21").Value = -0.0090042625818744781
$ws.Range("B21").Value = 0.0089999991833966675
$ws.Range("A22").Value = -0.15751647554113468
$ws.Range("B22").Value = 0.15645327012074262
$ws.Range("A23").Value = -0.0846426249356238
$ws.Range("B23").Value = 0.084128258459808869
$ws.Range("A24").Value = -0.042128259636220911
$ws.Range("B24").Value = 0.041999998817220607
$ws.Range("A25").Value = -0.05678455648548919
$ws.Range("B25").Value = 0.05669648115844339
$ws.Range("A26").Value = -0.050696481998475207
$ws.Range("B26").Value = 0.050586873354742323
$ws.Range("A27").Value = -0.044586874196887116
$ws.Range("B27").Value = 0.044222730725610671
$ws.Range("A28").Value = -0.038222731577056912
$ws.Range("B28").Value = 0.037985755470750782
$ws.Range("A29").Value = -0.02598575638642231
$ws.Range("B29").Value = 0.025889321023582212
$ws.Range("A30").Value = -0.042168393465581211
$ws.Range("B30").Value = 0.04201875407755562
$ws.Range("A31").Value = -0.027018755032777619
$ws.Range("B31").Value = 0.0270005036653167
$ws.Range("A32").Value = -0.0060005046800810646
$ws.Range("B32").Value = 0.0059999991308510303
